# Append the latest weekly bitcoin buy entry as a new row at the bottom
# of the data table (mirrors the existing rows' layout: Date, Coins, Price, Cost).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (row 70 -> 71).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Column A: the date is stored as a literal text string (like the other
# recent rows), not an actual date serial value, so force Text formatting
# before assigning it.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "02/22/2026"

# Column B: number of coins bought this week.
$ws.Cells.Item($newRow, 2).Value = 0.0007252900000000034

# Column C: price per coin at time of purchase.
$ws.Cells.Item($newRow, 3).Value = 68248.56264390763

# Column D: total cost of the purchase (USD).
$ws.Cells.Item($newRow, 4).Value = 50
